$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H41").Value = 160.625
$ws.Range("I41").Value = 144.5
$ws.Range("J41").Value = 166
$ws.Range("K41").Value = 144.5
$ws.Range("L41").Value = 166
$ws.Range("M41").Value = 295.5
$ws.Range("N41").Value = -1046

$ws.Range("H42").Value = 433.33334
$ws.Range("I42").Value = 400
$ws.Range("K42").Value = 1200
$ws.Range("M42").Value = -970

$ws.Range("H76").Value = 3006.4285
$ws.Range("I76").Value = 3199.25
$ws.Range("K76").Value = 3199.25
$ws.Range("M76").Value = -2884.25

$ws.Range("H79").Value = 3006.4285
$ws.Range("I79").Value = 3199.25
$ws.Range("K79").Value = 3199.25
$ws.Range("M79").Value = -2107.25

$ws.Range("H80").Value = 735.6667
$ws.Range("I80").Value = 421.2857
$ws.Range("J80").Value = 1175.8
$ws.Range("K80").Value = 1263.8571
$ws.Range("L80").Value = 3527.4
$ws.Range("M80").Value = -265.8571000000002
$ws.Range("N80").Value = -5523.4

$ws.Range("H83").Value = 735.6667
$ws.Range("I83").Value = 421.2857
$ws.Range("J83").Value = 1175.8
$ws.Range("K83").Value = 3791.5713
$ws.Range("L83").Value = 10582.2
$ws.Range("M83").Value = 1200.4287
$ws.Range("N83").Value = -20566.2

$ws.Range("H98").Value = 4990.0527
$ws.Range("I98").Value = 4552.6875
$ws.Range("K98").Value = 4552.6875
$ws.Range("M98").Value = -3054.6875

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("M106").ClearContents()
$ws.Range("N106").ClearContents()

$ws.Range("H107").Value = 296.33334
$ws.Range("I107").Value = 296.33334
$ws.Range("K107").Value = 296.33334
$ws.Range("M107").Value = 1623.66666

$ws.Range("H122").Value = 4990.0527
$ws.Range("I122").Value = 4552.6875
$ws.Range("K122").Value = 13658.0625
$ws.Range("M122").Value = -11208.0625

$ws.Range("H127").Value = 538.625
$ws.Range("I127").Value = 330
$ws.Range("J127").Value = 1999
$ws.Range("K127").Value = 990
$ws.Range("L127").Value = 5997
$ws.Range("M127").Value = 3970
$ws.Range("N127").Value = -15917

$ws.Range("H131").Value = 564012.7
$ws.Range("I131").Value = 722766.0600000001
$ws.Range("K131").Value = 2168298.18
$ws.Range("M131").Value = -2163258.18

$ws.Range("H132").Value = 4741.879
$ws.Range("I132").Value = 4858.8125
$ws.Range("K132").Value = 14576.4375
$ws.Range("M132").Value = -12046.4375

$ws.Range("H137").Value = 1434139
$ws.Range("I137").Value = 1925004
$ws.Range("K137").Value = 5775012
$ws.Range("M137").Value = -5772462

$ws.Range("H138").Value = 2852.6326
$ws.Range("I138").Value = 1399
$ws.Range("K138").Value = 4197
$ws.Range("M138").Value = 943

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 147976.2
$ws.Range("I74").Value = 151897.19
$ws.Range("K74").Value = 151897.19
$ws.Range("M74").Value = -151023.19

$ws.Range("H77").Value = 147976.2
$ws.Range("I77").Value = 151897.19
$ws.Range("K77").Value = 759485.95
$ws.Range("M77").Value = -755117.95

$ws.Range("H97").Value = 714.7857
$ws.Range("J97").Value = 850
$ws.Range("L97").Value = 850
$ws.Range("N97").Value = -1842

$ws.Range("H102").Value = 2874.6155
$ws.Range("I102").Value = 2442.7273
$ws.Range("K102").Value = 2442.7273
$ws.Range("M102").Value = -820.7273

$ws.Range("H110").Value = 6019.1
$ws.Range("I110").Value = 5575.3335
$ws.Range("K110").Value = 5575.3335
$ws.Range("M110").Value = -3530.3335

$ws.Range("H122").Value = 2357.4285
$ws.Range("I122").Value = 2385.261
$ws.Range("K122").Value = 7155.782999999999
$ws.Range("M122").Value = -4705.782999999999

$ws.Range("H132").Value = 2372.5
$ws.Range("I132").Value = 1528.3334
$ws.Range("J132").Value = 3216.6667
$ws.Range("K132").Value = 4585.0002
$ws.Range("L132").Value = 9650.000100000001
$ws.Range("M132").Value = -2055.0002
$ws.Range("N132").Value = -14710.0001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("K15").Value = 0
$ws.Range("M15").ClearContents()

$ws.Range("H94").Value = 142865310
$ws.Range("I94").Value = 222233980
$ws.Range("K94").Value = 222233980
$ws.Range("M94").Value = -222233529

$ws.Range("H99").Value = 7798.3335
$ws.Range("I99").Value = 7947.5
$ws.Range("K99").Value = 7947.5
$ws.Range("M99").Value = -6449.5

$ws.Range("H105").Value = 13001539
$ws.Range("I105").Value = 834519.5600000001
$ws.Range("J105").Value = 31252068
$ws.Range("K105").Value = 834519.5600000001
$ws.Range("L105").Value = 31252068
$ws.Range("M105").Value = -832772.5600000001
$ws.Range("N105").Value = -31255562

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1842.75
$ws.Range("J16").Value = 1937.5
$ws.Range("L16").Value = 1937.5
$ws.Range("N16").Value = -2511.5

$ws.Range("H31").Value = 4109.7114
$ws.Range("I31").Value = 3118.6072
$ws.Range("J31").Value = 5266
$ws.Range("K31").Value = 3118.6072
$ws.Range("L31").Value = 5266
$ws.Range("M31").Value = -2823.6072
$ws.Range("N31").Value = -5856

$ws.Range("H34").Value = 4109.7114
$ws.Range("I34").Value = 3118.6072
$ws.Range("J34").Value = 5266
$ws.Range("K34").Value = 3118.6072
$ws.Range("L34").Value = 5266
$ws.Range("M34").Value = -2916.6072
$ws.Range("N34").Value = -5670

$ws.Range("H105").Value = 2740.2222
$ws.Range("I105").Value = 2728
$ws.Range("K105").Value = 2728
$ws.Range("M105").Value = -981

$ws.Range("H113").Value = 1842.75
$ws.Range("J113").Value = 1937.5
$ws.Range("L113").Value = 1937.5
$ws.Range("N113").Value = -6277.5

$ws.Range("H122").Value = 2147.8333
$ws.Range("I122").Value = 2200.9
$ws.Range("J122").Value = 1882.5
$ws.Range("K122").Value = 6602.700000000001
$ws.Range("L122").Value = 5647.5
$ws.Range("M122").Value = -4152.700000000001
$ws.Range("N122").Value = -10547.5

$ws.Range("H132").Value = 14498380
$ws.Range("I132").Value = 4866.0835
$ws.Range("K132").Value = 14598.2505
$ws.Range("M132").Value = -12068.2505

$ws.Range("H134").Value = 4359.316
$ws.Range("I134").Value = 4762.2
$ws.Range("K134").Value = 14286.6
$ws.Range("M134").Value = -11751.6

$ws.Range("H139").Value = 48999.5
$ws.Range("J139").Value = 48999.5
$ws.Range("L139").Value = 48999.5
$ws.Range("N139").Value = -59279.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 5162.5
$ws.Range("J42").Value = 5162.5
$ws.Range("L42").Value = 15487.5
$ws.Range("N42").Value = -16555.5

$ws.Range("H121").Value = 4213399.5
$ws.Range("J121").Value = 69675.94
$ws.Range("L121").Value = 209027.82
$ws.Range("N121").Value = -211647.82

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 760.5
$ws.Range("I97").Value = 818.5714
$ws.Range("J97").Value = 354
$ws.Range("K97").Value = 818.5714
$ws.Range("L97").Value = 354
$ws.Range("M97").Value = -322.5714
$ws.Range("N97").Value = -1346

$ws.Range("H113").Value = 2605
$ws.Range("I113").Value = 2499.5
$ws.Range("K113").Value = 2499.5
$ws.Range("M113").Value = -329.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2326.125
$ws.Range("I7").Value = 1934.8334
$ws.Range("K7").Value = 1934.8334
$ws.Range("M7").Value = -1822.8334

$ws.Range("H22").Value = 1712.5
$ws.Range("I22").Value = 2450
$ws.Range("K22").Value = 2450
$ws.Range("M22").Value = -2155

$ws.Range("H27").Value = 1712.5
$ws.Range("I27").Value = 2450
$ws.Range("K27").Value = 2450
$ws.Range("M27").Value = -2343

$ws.Range("H40").Value = 25001.334
$ws.Range("I40").Value = 36002
$ws.Range("J40").Value = 3000
$ws.Range("K40").Value = 36002
$ws.Range("L40").Value = 3000
$ws.Range("M40").Value = -35866
$ws.Range("N40").Value = -3272

$ws.Range("H61").Value = 4148.077
$ws.Range("I61").Value = 4211.364
$ws.Range("K61").Value = 4211.364
$ws.Range("M61").Value = -4009.364

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H113").Value = 4148.077
$ws.Range("I113").Value = 4211.364
$ws.Range("K113").Value = 4211.364
$ws.Range("M113").Value = -2041.364

$ws.Range("H126").Value = 2326.125
$ws.Range("I126").Value = 1934.8334
$ws.Range("K126").Value = 5804.5002
$ws.Range("M126").Value = -3334.5002

$ws.Range("H136").Value = 5168.2173
$ws.Range("I136").Value = 4127.636
$ws.Range("K136").Value = 12382.908
$ws.Range("M136").Value = -9832.908000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 17755.143
$ws.Range("J41").Value = 17303.6
$ws.Range("L41").Value = 17303.6
$ws.Range("N41").Value = -18083.6

$ws.Range("H122").Value = 8336413
$ws.Range("I122").Value = 3264.577
$ws.Range("K122").Value = 9793.731
$ws.Range("M122").Value = -7343.731
